$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cell from "教科" to "教科名"
$ws.Range("A1").Value = "教科名"

# Move the active selection to F10, matching the author's last cursor position
[void]$ws.Range("F10").Select()
